$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.855.78'
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '2.179.72'
$ws.Range("E3").Value = '  -2.80%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'238.44"
$ws.Range("E5").Value = '  -1.89%  '

$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = '  -1.94%  '

$ws.Range("D7").Value = "'72.81"
$ws.Range("E7").Value = '  -1.84%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = '  -2.92%  '

$ws.Range("D10").Value = "'40.02"
$ws.Range("E10").Value = '  -5.61%  '

$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = '  -4.98%  '

$ws.Range("D12").Value = "'54.64"
$ws.Range("E12").Value = '  -3.39%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'6.73"
$ws.Range("E13").Value = '  -2.99%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = "'0.100"
$ws.Range("E14").Value = '  -3.32%  '

$ws.Range("D15").Value = '2.501.84'
$ws.Range("E15").Value = '  -3.00%  '

$ws.Range("D16").Value = "'14.38"
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '2.174.15'
$ws.Range("E17").Value = '  -3.10%  '

$ws.Range("D18").Value = "'0.784"
$ws.Range("E18").Value = '  -6.73%  '

$ws.Range("D19").Value = '41.675.81'
$ws.Range("E19").Value = '  -0.95%  '

$ws.Range("D20").Value = "'0.0000103"
$ws.Range("E20").Value = '  -1.79%  '

$ws.Range("D21").Value = "'70.12"
$ws.Range("E21").Value = '  -3.59%  '

$ws.Range("D22").Value = "'5.81"
$ws.Range("E22").Value = '  -6.76%  '

$ws.Range("E23").Value = '  -11.25%  '

$ws.Range("D24").Value = "'226.45"
$ws.Range("E24").Value = '  -1.71%  '

$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = '  -6.01%  '

$ws.Range("D28").Value = "'3.28"
$ws.Range("E28").Value = '  -9.39%  '

$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = '  -3.58%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'171.74"
$ws.Range("E30").Value = '  +2.67%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'2.08"
$ws.Range("E31").Value = '  -5.60%  '

$ws.Range("D32").Value = "'19.92"
$ws.Range("E32").Value = '  -3.37%  '

$ws.Range("D33").Value = "'32.88"
$ws.Range("E33").Value = '  +10.66%  '

$ws.Range("D34").Value = "'0.0776"
$ws.Range("E34").Value = '  -3.63%  '

$ws.Range("D35").Value = "'5.28"
$ws.Range("E35").Value = '  -6.23%  '

$ws.Range("D36").Value = "'0.121"
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("D37").Value = "'4.31"
$ws.Range("E37").Value = '  -0.50%  '

$ws.Range("E38").Value = '  -6.63%  '

$ws.Range("D39").Value = "'0.0312"
$ws.Range("E39").Value = '  +2.37%  '

$ws.Range("D40").Value = "'12.12"
$ws.Range("E40").Value = '  -8.29%  '

$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("D42").Value = "'5.37"
$ws.Range("E42").Value = '  -6.08%  '

$ws.Range("D43").Value = "'59.20"
$ws.Range("E43").Value = '  -8.66%  '

$ws.Range("D44").Value = "'0.190"
$ws.Range("E44").Value = '  -4.57%  '

$ws.Range("D45").Value = "'8.43"
$ws.Range("E45").Value = '  -3.55%  '

$ws.Range("D46").Value = "'0.0966"
$ws.Range("E46").Value = '  -3.91%  '

$ws.Range("D47").Value = "'97.49"
$ws.Range("E47").Value = '  -6.55%  '

$ws.Range("E48").Value = '  -4.60%  '

$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = '  -4.81%  '

$ws.Range("D50").Value = "'2.21"
$ws.Range("E50").Value = '  -5.38%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.412"
$ws.Range("E51").Value = '  +11.33%  '
